$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" values (regenerated from source data to replace the old "Strike#" column)
$kValues = @{
    2 = 1
    3 = 2
    4 = 4
    5 = 2
    6 = 0
    7 = 1
    8 = 1
    9 = 1
    10 = 3
    11 = 2
    12 = 1
    13 = 0
    14 = 0
    15 = 1
    16 = 3
    17 = 1
    18 = 2
    19 = 0
    20 = 1
    21 = 2
    22 = 0
    23 = 0
    24 = 1
    25 = 2
    26 = 1
    27 = 2
    28 = 0
    29 = 2
    30 = 2
    31 = 0
    32 = 0
    33 = 1
    34 = 1
    35 = 0
    36 = 1
    37 = 2
    38 = 2
    39 = 3
    40 = 1
    41 = 2
    42 = 3
    43 = 0
    44 = 3
    45 = 1
    46 = 0
    47 = 1
    48 = 0
    49 = 2
    50 = 0
    51 = 0
    52 = 1
    53 = 1
    54 = 0
    55 = 1
    56 = 1
    57 = 1
    58 = 0
    59 = 2
    60 = 0
    61 = 0
    62 = 0
    63 = 0
    64 = 0
    65 = 0
    66 = 2
    67 = 2
    68 = 1
    69 = 0
    70 = 1
    71 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
